# Auto-generated script to update FFXIV Leve profit market-data cells
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 443.66666
$ws.Range("I32").Value = 425.25
$ws.Range("J32").Value = 458.4
$ws.Range("K32").Value = 425.25
$ws.Range("L32").Value = 458.4
$ws.Range("M32").Value = -99.25
$ws.Range("N32").Value = -1110.4
$ws.Range("H51").Value = 7796
$ws.Range("I51").Value = 1980
$ws.Range("J51").Value = 9250
$ws.Range("K51").Value = 1980
$ws.Range("L51").Value = 9250
$ws.Range("M51").Value = -1496
$ws.Range("N51").Value = -10218
$ws.Range("H132").Value = 20491912
$ws.Range("I132").Value = 23350486
$ws.Range("J132").Value = 5463.3335
$ws.Range("K132").Value = 70051458
$ws.Range("L132").Value = 16390.0005
$ws.Range("M132").Value = -70048928
$ws.Range("N132").Value = -21450.0005
$ws.Range("H138").Value = 2445.404
$ws.Range("I138").Value = 1298.08
$ws.Range("J138").Value = 2833.0134
$ws.Range("K138").Value = 3894.24
$ws.Range("L138").Value = 8499.040199999999
$ws.Range("M138").Value = 1245.76
$ws.Range("N138").Value = -18779.0402

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2127.6155
$ws.Range("I74").Value = 1469.9546
$ws.Range("J74").Value = 5744.75
$ws.Range("K74").Value = 1469.9546
$ws.Range("L74").Value = 5744.75
$ws.Range("M74").Value = -595.9546
$ws.Range("N74").Value = -7492.75
$ws.Range("H77").Value = 2127.6155
$ws.Range("I77").Value = 1469.9546
$ws.Range("J77").Value = 5744.75
$ws.Range("K77").Value = 7349.773
$ws.Range("L77").Value = 28723.75
$ws.Range("M77").Value = -2981.773
$ws.Range("N77").Value = -37459.75
$ws.Range("H122").Value = 2162.7917
$ws.Range("I122").Value = 1349.8948
$ws.Range("J122").Value = 5251.8
$ws.Range("K122").Value = 4049.6844
$ws.Range("L122").Value = 15755.4
$ws.Range("M122").Value = -1599.6844
$ws.Range("N122").Value = -20655.4
$ws.Range("H132").Value = 2288.3257
$ws.Range("I132").Value = 1165.4138
$ws.Range("J132").Value = 4614.357
$ws.Range("K132").Value = 3496.2414
$ws.Range("L132").Value = 13843.071
$ws.Range("M132").Value = -966.2413999999999
$ws.Range("N132").Value = -18903.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1822.6111
$ws.Range("I105").Value = 1708.4286
$ws.Range("J105").Value = 2222.25
$ws.Range("K105").Value = 1708.4286
$ws.Range("L105").Value = 2222.25
$ws.Range("M105").Value = 38.57140000000004
$ws.Range("N105").Value = -5716.25
$ws.Range("H132").Value = 42954.547
$ws.Range("J132").Value = 42954.547
$ws.Range("L132").Value = 42954.547
$ws.Range("N132").Value = -53074.547
$ws.Range("H134").Value = 2647.102
$ws.Range("I134").Value = 1458.8206
$ws.Range("J134").Value = 7281.4
$ws.Range("K134").Value = 4376.4618
$ws.Range("L134").Value = 21844.2
$ws.Range("M134").Value = -1841.4618
$ws.Range("N134").Value = -26914.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 460.84616
$ws.Range("I22").Value = 309.12
$ws.Range("J22").Value = 731.7857
$ws.Range("K22").Value = 309.12
$ws.Range("L22").Value = 731.7857
$ws.Range("M22").Value = 40.88
$ws.Range("N22").Value = -1431.7857
$ws.Range("H31").Value = 2768.492
$ws.Range("I31").Value = 1290.8158
$ws.Range("J31").Value = 5014.56
$ws.Range("K31").Value = 1290.8158
$ws.Range("L31").Value = 5014.56
$ws.Range("M31").Value = -995.8158000000001
$ws.Range("N31").Value = -5604.56
$ws.Range("H34").Value = 2768.492
$ws.Range("I34").Value = 1290.8158
$ws.Range("J34").Value = 5014.56
$ws.Range("K34").Value = 1290.8158
$ws.Range("L34").Value = 5014.56
$ws.Range("M34").Value = -1088.8158
$ws.Range("N34").Value = -5418.56
$ws.Range("H58").Value = 1833.5375
$ws.Range("I58").Value = 1618.1471
$ws.Range("J58").Value = 3054.0833
$ws.Range("K58").Value = 1618.1471
$ws.Range("L58").Value = 3054.0833
$ws.Range("M58").Value = -1415.1471
$ws.Range("N58").Value = -3460.0833
$ws.Range("H100").Value = 62000
$ws.Range("J100").Value = 62000
$ws.Range("L100").Value = 62000
$ws.Range("N100").Value = -64164
$ws.Range("H122").Value = 2420.68
$ws.Range("I122").Value = 1808.6471
$ws.Range("K122").Value = 5425.9413
$ws.Range("M122").Value = -2975.9413
$ws.Range("H132").Value = 3446
$ws.Range("I132").Value = 2840
$ws.Range("K132").Value = 8520
$ws.Range("M132").Value = -5990
$ws.Range("H134").Value = 4682.2646
$ws.Range("I134").Value = 4998.923
$ws.Range("J134").Value = 3653.125
$ws.Range("K134").Value = 14996.769
$ws.Range("L134").Value = 10959.375
$ws.Range("M134").Value = -12461.769
$ws.Range("N134").Value = -16029.375
$ws.Range("H136").Value = 1833.5375
$ws.Range("I136").Value = 1618.1471
$ws.Range("J136").Value = 3054.0833
$ws.Range("K136").Value = 4854.4413
$ws.Range("L136").Value = 9162.249899999999
$ws.Range("M136").Value = -2304.4413
$ws.Range("N136").Value = -14262.2499
$ws.Range("H139").Value = 49380
$ws.Range("J139").Value = 49380
$ws.Range("L139").Value = 49380
$ws.Range("N139").Value = -59660

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 51.5
$ws.Range("J12").Value = 30.5
$ws.Range("L12").Value = 91.5
$ws.Range("N12").Value = -437.5
$ws.Range("H113").Value = 546.5806
$ws.Range("I113").Value = 564.41174
$ws.Range("J113").Value = 524.9286
$ws.Range("K113").Value = 1693.23522
$ws.Range("L113").Value = 1574.7858
$ws.Range("M113").Value = 476.76478
$ws.Range("N113").Value = -5914.7858
$ws.Range("H131").Value = 10870573
$ws.Range("I131").Value = 29412846
$ws.Range("J131").Value = 965.0345
$ws.Range("K131").Value = 88238538
$ws.Range("L131").Value = 2895.1035
$ws.Range("M131").Value = -88233498
$ws.Range("N131").Value = -12975.1035
$ws.Range("H134").Value = 2901.0605
$ws.Range("I134").Value = 1910.6818
$ws.Range("J134").Value = 4881.8184
$ws.Range("K134").Value = 5732.0454
$ws.Range("L134").Value = 14645.4552
$ws.Range("M134").Value = -662.0454
$ws.Range("N134").Value = -24785.4552
$ws.Range("H140").Value = 21418.04
$ws.Range("I140").Value = 42905.75
$ws.Range("K140").Value = 128717.25
$ws.Range("M140").Value = -123537.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 22200.416
$ws.Range("J43").Value = 25570.5
$ws.Range("L43").Value = 25570.5
$ws.Range("N43").Value = -25872.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1639.5758
$ws.Range("I22").Value = 1321.4117
$ws.Range("J22").Value = 1977.625
$ws.Range("K22").Value = 1321.4117
$ws.Range("L22").Value = 1977.625
$ws.Range("M22").Value = -1026.4117
$ws.Range("N22").Value = -2567.625
$ws.Range("H27").Value = 1639.5758
$ws.Range("I27").Value = 1321.4117
$ws.Range("J27").Value = 1977.625
$ws.Range("K27").Value = 1321.4117
$ws.Range("L27").Value = 1977.625
$ws.Range("M27").Value = -1214.4117
$ws.Range("N27").Value = -2191.625
$ws.Range("H122").Value = 3642.3635
$ws.Range("I122").Value = 3138.6453
$ws.Range("J122").Value = 11450
$ws.Range("K122").Value = 9415.9359
$ws.Range("L122").Value = 34350
$ws.Range("M122").Value = -6965.9359
$ws.Range("N122").Value = -39250
$ws.Range("H132").Value = 4764.381
$ws.Range("I132").Value = 1954.8
$ws.Range("J132").Value = 7318.5454
$ws.Range("K132").Value = 5864.4
$ws.Range("L132").Value = 21955.6362
$ws.Range("M132").Value = -3334.4
$ws.Range("N132").Value = -27015.6362
$ws.Range("H136").Value = 3446.2432
$ws.Range("I136").Value = 1931.8948
$ws.Range("J136").Value = 5044.722
$ws.Range("K136").Value = 5795.6844
$ws.Range("L136").Value = 15134.166
$ws.Range("M136").Value = -3245.6844
$ws.Range("N136").Value = -20234.166

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 82423.89
$ws.Range("J46").Value = 82423.89
$ws.Range("L46").Value = 82423.89
$ws.Range("N46").Value = -82885.89
$ws.Range("H107").Value = 1058.3572
$ws.Range("I107").Value = 991.1
$ws.Range("K107").Value = 2973.3
$ws.Range("M107").Value = -1053.3
$ws.Range("H122").Value = 6326.4736
$ws.Range("I122").Value = 5013.533
$ws.Range("J122").Value = 11250
$ws.Range("K122").Value = 15040.599
$ws.Range("L122").Value = 33750
$ws.Range("M122").Value = -12590.599
$ws.Range("N122").Value = -38650
$ws.Range("H123").Value = 29935
$ws.Range("J123").Value = 29935
$ws.Range("L123").Value = 29935
$ws.Range("N123").Value = -39735
$ws.Range("H126").Value = 329233.2
$ws.Range("I126").Value = 2953
$ws.Range("J126").Value = 601133.3
$ws.Range("K126").Value = 8859
$ws.Range("L126").Value = 1803399.9
$ws.Range("M126").Value = -6389
$ws.Range("N126").Value = -1808339.9
$ws.Range("H134").Value = 82423.89
$ws.Range("J134").Value = 82423.89
$ws.Range("L134").Value = 247271.67
$ws.Range("N134").Value = -252341.67
$ws.Range("H136").Value = 3740.8518
$ws.Range("I136").Value = 1621.2307
$ws.Range("J136").Value = 5709.0713
$ws.Range("K136").Value = 4863.6921
$ws.Range("L136").Value = 17127.2139
$ws.Range("M136").Value = -2313.6921
$ws.Range("N136").Value = -22227.2139
